$wb = $excel.ActiveWorkbook

# ---- Overview ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Rows(5).Insert()

$ws.Range("A5").Value = "13cc79f3-337b-4661-a926-917e426cbb9d.md"
$ws.Range("B5").Value = "e2e\13cc79f3-337b-4661-a926-917e426cbb9d.md"
$ws.Range("C5").Value = ".md"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-20 10:47:32"

$ws.Range("A6").Value = "4cf8f5d3-6d3e-4436-9b05-b7100652348b.md"
$ws.Range("B6").Value = "e2e\4cf8f5d3-6d3e-4436-9b05-b7100652348b.md"
$ws.Range("C6").Value = ".md"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "Ready for handoff"
$ws.Range("F6").Value = "Ready for handoff"
$ws.Range("G6").Value = "2016-08-20 10:45:54"

$ws.Range("A7").Value = "54d82c55-4b6e-431e-97bd-4f52a4772ec0.md"
$ws.Range("B7").Value = "e2e\54d82c55-4b6e-431e-97bd-4f52a4772ec0.md"
$ws.Range("C7").Value = ".md"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "Ready for handoff"
$ws.Range("F7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2016-08-20 10:47:32"

$t = $ws.ListObjects.Item("Overview")
$t.Resize($ws.Range("A1:G7"))

# ---- zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Rows(5).Insert()

$ws.Range("A5").Value = "13cc79f3-337b-4661-a926-917e426cbb9d.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "e2e"
$ws.Range("E5").Value = "ht"
$ws.Range("F5").Value = "False"
$ws.Range("G5").Value = "13cc79f3-337b-4661-a926-917e426cbb9d.01b8d2cb9de76d8ecfab78ea5a7c520a79e09aad.zh-cn.xlf"
$ws.Range("H5").Value = "2016-08-20 10:47:28"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "0001-01-01 00:00:00"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "True"
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "False"
$ws.Range("P5").Value = ""

$ws.Range("A6").Value = "4cf8f5d3-6d3e-4436-9b05-b7100652348b.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "e2e"
$ws.Range("E6").Value = "ht"
$ws.Range("F6").Value = "False"
$ws.Range("G6").Value = "4cf8f5d3-6d3e-4436-9b05-b7100652348b.157e76ccf22cb4ac9cfa329bd62fd836d8e1fa46.zh-cn.xlf"
$ws.Range("H6").Value = "2016-08-20 10:45:50"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = "0001-01-01 00:00:00"
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = "True"
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = "False"
$ws.Range("P6").Value = ""

$ws.Range("A7").Value = "54d82c55-4b6e-431e-97bd-4f52a4772ec0.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "e2e"
$ws.Range("E7").Value = "ht"
$ws.Range("F7").Value = "False"
$ws.Range("G7").Value = "54d82c55-4b6e-431e-97bd-4f52a4772ec0.79241a373340ada82b09d0d098f7aa52c985f315.zh-cn.xlf"
$ws.Range("H7").Value = "2016-08-20 10:47:28"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = "0001-01-01 00:00:00"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = "True"
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = "False"
$ws.Range("P7").Value = ""

$t = $ws.ListObjects.Item("zh-cn")
$t.Resize($ws.Range("A1:P7"))

# ---- de-de ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Rows(5).Insert()

$ws.Range("A5").Value = "13cc79f3-337b-4661-a926-917e426cbb9d.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "e2e"
$ws.Range("E5").Value = "ht"
$ws.Range("F5").Value = "False"
$ws.Range("G5").Value = "13cc79f3-337b-4661-a926-917e426cbb9d.01b8d2cb9de76d8ecfab78ea5a7c520a79e09aad.de-de.xlf"
$ws.Range("H5").Value = "2016-08-20 10:47:32"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "0001-01-01 00:00:00"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "True"
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "False"
$ws.Range("P5").Value = ""

$ws.Range("A6").Value = "4cf8f5d3-6d3e-4436-9b05-b7100652348b.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "e2e"
$ws.Range("E6").Value = "ht"
$ws.Range("F6").Value = "False"
$ws.Range("G6").Value = "4cf8f5d3-6d3e-4436-9b05-b7100652348b.157e76ccf22cb4ac9cfa329bd62fd836d8e1fa46.de-de.xlf"
$ws.Range("H6").Value = "2016-08-20 10:45:54"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = "0001-01-01 00:00:00"
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = "True"
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = "False"
$ws.Range("P6").Value = ""

$ws.Range("A7").Value = "54d82c55-4b6e-431e-97bd-4f52a4772ec0.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "e2e"
$ws.Range("E7").Value = "ht"
$ws.Range("F7").Value = "False"
$ws.Range("G7").Value = "54d82c55-4b6e-431e-97bd-4f52a4772ec0.79241a373340ada82b09d0d098f7aa52c985f315.de-de.xlf"
$ws.Range("H7").Value = "2016-08-20 10:47:32"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = "0001-01-01 00:00:00"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = "True"
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = "False"
$ws.Range("P7").Value = ""

$t = $ws.ListObjects.Item("de-de")
$t.Resize($ws.Range("A1:P7"))
